$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking strings are preserved as text (matching source inlineStr cells)
$valueRange = $ws.Range("D2:E51")
$valueRange.NumberFormat = "@"

$ws.Range("D2").Value = '66.847.42'
$ws.Range("E2").Value = '  -4.73%  '
$ws.Range("D3").Value = '3.209.78'
$ws.Range("E3").Value = '  -8.53%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '591.34'
$ws.Range("E5").Value = '  -2.36%  '
$ws.Range("D6").Value = '150.37'
$ws.Range("E6").Value = '  -13.01%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '3.199.56'
$ws.Range("E8").Value = '  -8.70%  '
$ws.Range("D9").Value = '0.542'
$ws.Range("E9").Value = '  -10.87%  '
$ws.Range("D10").Value = '0.172'
$ws.Range("E10").Value = '  -11.66%  '
$ws.Range("D11").Value = '6.50'
$ws.Range("E11").Value = '  -10.20%  '
$ws.Range("D12").Value = '0.498'
$ws.Range("E12").Value = '  -15.05%  '
$ws.Range("D13").Value = '38.88'
$ws.Range("E13").Value = '  -15.77%  '
$ws.Range("D14").Value = '0.0000243'
$ws.Range("E14").Value = '  -11.79%  '
$ws.Range("D15").Value = '3.733.70'
$ws.Range("E15").Value = '  -8.45%  '
$ws.Range("D16").Value = '66.895.61'
$ws.Range("E16").Value = '  -4.67%  '
$ws.Range("D17").Value = '3.213.22'
$ws.Range("E17").Value = '  -8.42%  '
$ws.Range("E18").Value = '  -5.09%  '
$ws.Range("D19").Value = '7.18'
$ws.Range("E19").Value = '  -14.18%  '
$ws.Range("D20").Value = '526.78'
$ws.Range("E20").Value = '  -14.15%  '
$ws.Range("D21").Value = '14.98'
$ws.Range("E21").Value = '  -14.35%  '
$ws.Range("D22").Value = '0.758'
$ws.Range("E22").Value = '  -13.70%  '
$ws.Range("D23").Value = '7.90'
$ws.Range("E23").Value = '  -13.79%  '
$ws.Range("D24").Value = '13.75'
$ws.Range("E24").Value = '  -11.68%  '
$ws.Range("D25").Value = '85.37'
$ws.Range("E25").Value = '  -13.38%  '
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("D27").Value = '3.17'
$ws.Range("E27").Value = '  -14.66%  '
$ws.Range("D28").Value = '2.21'
$ws.Range("E28").Value = '  -13.73%  '
$ws.Range("D29").Value = '8.07'
$ws.Range("E29").Value = '  -10.19%  '
$ws.Range("D30").Value = '29.23'
$ws.Range("E30").Value = '  -13.71%  '
$ws.Range("D31").Value = '2.65'
$ws.Range("E31").Value = '  -11.14%  '
$ws.Range("D32").Value = '1.13'
$ws.Range("E32").Value = '  -11.19%  '
$ws.Range("D33").Value = '545.08'
$ws.Range("E33").Value = '  -13.84%  '
$ws.Range("D34").Value = '6.50'
$ws.Range("E34").Value = '  -19.34%  '
$ws.Range("D35").Value = '5.69'
$ws.Range("E35").Value = '  -16.45%  '
$ws.Range("E36").Value = '  +0.33%  '
$ws.Range("D37").Value = '53.34'
$ws.Range("E37").Value = '  -6.07%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Value = '0.0863'
$ws.Range("E38").Value = '  -13.52%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.0423'
$ws.Range("E39").Value = '  -12.23%  '
$ws.Range("D40").Value = '9.32'
$ws.Range("E40").Value = '  -13.34%  '
$ws.Range("D41").Value = '0.128'
$ws.Range("E41").Value = '  -11.87%  '
$ws.Range("D42").Value = '2.923.33'
$ws.Range("E42").Value = '  -13.10%  '
$ws.Range("D43").Value = '2.62'
$ws.Range("E43").Value = '  -24.64%  '
$ws.Range("E44").Value = '  -15.24%  '
$ws.Range("D45").Value = '0.0₃0582'
$ws.Range("E45").Value = '  -20.70%  '
$ws.Range("D46").Value = '2.40'
$ws.Range("E46").Value = '  -17.56%  '
$ws.Range("D47").Value = '26.45'
$ws.Range("E47").Value = '  -17.15%  '
$ws.Range("E48").Value = '  -0.09%  '
$ws.Range("E49").Value = '  -17.43%  '
$ws.Range("E50").Value = '  -12.81%  '
$ws.Range("D51").Value = '123.46'
$ws.Range("E51").Value = '  -7.32%  '

# Restore default styling so no extraneous style index is left on the cells
$valueRange.Style = "Normal"
